$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 114, shifting existing rows 114:223 down to 115:224
$ws.Rows.Item(114).Insert()

# Populate the new row 114 with its data
$ws.Cells.Item(114, 1).Value2 = 8
$ws.Cells.Item(114, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(114, 3).Value2 = "Coquimbo"
$ws.Cells.Item(114, 4).Value2 = 45033
$ws.Cells.Item(114, 5).Value2 = 4
$ws.Cells.Item(114, 6).Value2 = 100112001
$ws.Cells.Item(114, 7).Value2 = "Berenjena"
$ws.Cells.Item(114, 8).Value2 = "Sin especificar"
$ws.Cells.Item(114, 9).Value2 = "Primera"
$ws.Cells.Item(114, 10).Value2 = 600
$ws.Cells.Item(114, 11).Value2 = 10000
$ws.Cells.Item(114, 12).Value2 = 11000
$ws.Cells.Item(114, 13).Value2 = 10500
$ws.Cells.Item(114, 14).Value2 = "$/caja 50 unidades"
$ws.Cells.Item(114, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(114, 16).Value2 = 210
$ws.Cells.Item(114, 17).Value2 = 50
$ws.Cells.Item(114, 18).Value2 = "Hortaliza"
